$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '27.692.39'
$c.ClearFormats()
$ws.Range("E2").Value = '  +0.02%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.901.64'
$c.ClearFormats()
$ws.Range("E3").Value = '  +0.78%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9996'
$c.ClearFormats()
$ws.Range("E4").Value = '  -0.17%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '312.36'
$c.ClearFormats()
$ws.Range("E5").Value = '  -0.25%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9996'
$c.ClearFormats()
$ws.Range("E6").Value = '  -0.16%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5225'
$c.ClearFormats()
$ws.Range("E7").Value = '  +7.92%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3779'
$c.ClearFormats()
$ws.Range("E8").Value = '  -0.22%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07245'
$c.ClearFormats()
$ws.Range("E9").Value = '  -1.06%  '
$ws.Range("E10").Value = '  +3.46%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.8952'
$c.ClearFormats()
$ws.Range("E11").Value = '  -2.53%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.926.16'
$c.ClearFormats()
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.07637'
$c.ClearFormats()
$ws.Range("E13").Value = '  -0.46%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.439'
$c.ClearFormats()
$ws.Range("E14").Value = '  -0.40%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '91.94'
$c.ClearFormats()
$ws.Range("E15").Value = '  +1.23%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.9998'
$c.ClearFormats()
$ws.Range("E16").Value = '  -0.25%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.000008716'
$c.ClearFormats()
$ws.Range("E17").Value = '  -0.88%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.9992'
$c.ClearFormats()
$ws.Range("E18").Value = '  -0.25%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '27.735.93'
$c.ClearFormats()
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  -0.39%  '
$ws.Range("E21").Value = '  +0.38%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '2.129.14'
$c.ClearFormats()
$ws.Range("E22").Value = '  -0.80%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '10.82'
$c.ClearFormats()
$ws.Range("E23").Value = '  +0.18%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '6.570'
$c.ClearFormats()
$ws.Range("E24").Value = '  -0.19%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '153.10'
$c.ClearFormats()
$ws.Range("E25").Value = '  -0.15%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.865'
$c.ClearFormats()
$ws.Range("E26").Value = '  -1.98%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.161'
$c.ClearFormats()
$ws.Range("E27").Value = '  +2.53%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '18.28'
$c.ClearFormats()
$ws.Range("E28").Value = '  -0.35%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '114.59'
$c.ClearFormats()
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '4.837'
$c.ClearFormats()
$ws.Range("E30").Value = '  -1.22%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.08990'
$c.ClearFormats()
$ws.Range("E31").Value = '  +0.66%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.170'
$c.ClearFormats()
$ws.Range("E32").Value = '  +0.77%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.240'
$c.ClearFormats()
$ws.Range("E33").Value = '  +1.65%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.808'
$c.ClearFormats()
$ws.Range("E34").Value = '  +4.08%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.7752'
$c.ClearFormats()
$ws.Range("E35").Value = '  +2.44%  '
$ws.Range("E36").Value = '  +2.42%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.607'
$c.ClearFormats()
$ws.Range("E37").Value = '  +2.76%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.055'
$c.ClearFormats()
$ws.Range("E38").Value = '  +2.87%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.090'
$c.ClearFormats()
$ws.Range("E39").Value = '  +0.11%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.5496'
$c.ClearFormats()
$ws.Range("E40").Value = '  +1.02%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.05284'
$c.ClearFormats()
$ws.Range("E41").Value = '  +0.81%  '
$ws.Range("E42").Value = '  -4.25%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '113.25'
$c.ClearFormats()
$ws.Range("E43").Value = '  +3.05%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '8.451'
$c.ClearFormats()
$ws.Range("E44").Value = '  +1.64%  '
$ws.Range("E45").Value = '  -0.64%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.4789'
$c.ClearFormats()
$ws.Range("E46").Value = '  +0.33%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '10.44'
$c.ClearFormats()
$ws.Range("E47").Value = '  -1.10%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.9990'
$c.ClearFormats()
$ws.Range("E48").Value = '  -0.19%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.614'
$c.ClearFormats()
$ws.Range("E49").Value = '  -0.70%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '66.47'
$c.ClearFormats()
$ws.Range("E50").Value = '  -1.46%  '
$ws.Range("E51").Value = '  -0.85%  '
